# Water Quality Monitor log: add entries for 2022-11-07, 2022-11-08, 2022-11-09
# and move the active selection to the newly filled-in E48 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Carry over the existing date/time number formats (and alignment) from the
# last populated row (46) so the new rows render the same way, instead of
# leaving the default "General" style that the blank placeholder rows have.
$ws.Range("B46:D46").Copy() | Out-Null
$ws.Range("B47:D49").PasteSpecial(-4122) | Out-Null

$ws.Range("E46").Copy() | Out-Null
$ws.Range("E47:E48").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# Row 47 - 2022-11-07, 12:00 - 17:00
$ws.Range("B47").Value = 44872
$ws.Range("C47").Value = 0.5
$ws.Range("D47").Value = 0.70833333333333337
$ws.Range("E47").Value = "Created test documents and did test for microscope led, made drawings. Fetched and created HCL boxes with components. Tried to fix pyqt5, again. "

# Row 48 - 2022-11-08, 09:15 - 15:00
$ws.Range("B48").Value = 44873
$ws.Range("C48").Value = 0.38541666666666669
$ws.Range("D48").Value = 0.625
$ws.Range("E48").Value = "fixed bug that let pigpiod to not behave properly. Wrote testreport on stepper motor and got stepper motor working very well."

# Row 49 - 2022-11-09, 09:30 - 17:00 (no notes entered yet)
$ws.Range("B49").Value = 44874
$ws.Range("C49").Value = 0.39583333333333331
$ws.Range("D49").Value = 0.70833333333333337

# Match the author's last selected cell.
$ws.Range("E48").Select() | Out-Null
